$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Range("D2").Value = '45.091.71'
    $ws.Range("E2").Value = '  +4.76%  '
    $ws.Range("D3").Value = '2.444.99'
    $ws.Range("E3").Value = '  +3.33%  '
    $ws.Range("D4").Value = '0.998'
    $ws.Range("E4").Value = '  -0.14%  '
    $ws.Range("D5").Value = '318.93'
    $ws.Range("E5").Value = '  +5.09%  '
    $ws.Range("D6").Value = '104.70'
    $ws.Range("E6").Value = '  +9.49%  '
    $ws.Range("E7").Value = '  +2.65%  '
    $ws.Range("E8").Value = '  -0.07%  '
    $ws.Range("D9").Value = '0.533'
    $ws.Range("E9").Value = '  +10.55%  '
    $ws.Range("D10").Value = '35.91'
    $ws.Range("E10").Value = '  +4.55%  '
    $ws.Range("E11").Value = '  +2.18%  '
    $ws.Range("E12").Value = '  -2.66%  '
    $ws.Range("E13").Value = '  +0.69%  '
    $ws.Range("D14").Value = '7.02'
    $ws.Range("E14").Value = '  +3.36%  '
    $ws.Range("D15").Value = '2.814.86'
    $ws.Range("E15").Value = '  +2.99%  '
    $ws.Range("D16").Value = '2.413.58'
    $ws.Range("E16").Value = '  +1.95%  '
    $ws.Range("D17").Value = '0.841'
    $ws.Range("E17").Value = '  +4.72%  '
    $ws.Range("D18").Value = '44.968.36'
    $ws.Range("E18").Value = '  +4.53%  '
    $ws.Range("D19").Value = '12.32'
    $ws.Range("E19").Value = '  +3.18%  '
    $ws.Range("D20").Value = '6.38'
    $ws.Range("E20").Value = '  +1.27%  '
    $ws.Range("E21").Value = '  +3.78%  '
    $ws.Range("D22").Value = '69.13'
    $ws.Range("D23").Value = '244.50'
    $ws.Range("E23").Value = '  +3.89%  '
    $ws.Range("D24").Value = '2.29'
    $ws.Range("E24").Value = '  +3.44%  '
    $ws.Range("D25").Value = '2.52'
    $ws.Range("E25").Value = '  +3.28%  '
    $ws.Range("E26").Value = '  +0.06%  '
    $ws.Range("D27").Value = '25.51'
    $ws.Range("E27").Value = '  +4.35%  '
    $ws.Range("D28").Value = '2.20'
    $ws.Range("E28").Value = '  -6.91%  '
    $ws.Range("D29").Value = '9.62'
    $ws.Range("E29").Value = '  +2.79%  '
    $ws.Range("D30").Value = '34.11'
    $ws.Range("E30").Value = '  +6.35%  '
    $ws.Range("D31").Value = '49.48'
    $ws.Range("E31").Value = '  +3.22%  '
    $ws.Range("D32").Value = '0.128'
    $ws.Range("E32").Value = '  +15.62%  '
    $ws.Range("D33").Value = '20.13'
    $ws.Range("E33").Value = '  +13.61%  '
    $ws.Range("E34").Value = '  +4.16%  '
    $ws.Range("E35").Value = '  +0.22%  '
    $ws.Range("E36").Value = '  +3.95%  '
    $ws.Range("E37").Value = '  +4.98%  '
    $ws.Range("E38").Value = '  +4.62%  '
    $ws.Range("E39").Value = '  +1.13%  '
    $ws.Range("D40").Value = '124.71'
    $ws.Range("E40").Value = '  -3.30%  '
    $ws.Range("E41").Value = '  +2.45%  '
    $ws.Range("D42").Value = '2.19'
    $ws.Range("E42").Value = '  -2.66%  '
    $ws.Range("D43").Value = '21.20'
    $ws.Range("E43").Value = '  +0.16%  '
    $ws.Range("E44").Value = '  +4.66%  '
    $ws.Range("D45").Value = '1.948.72'
    $ws.Range("E45").Value = '  +1.06%  '
    $ws.Range("E46").Value = '  +8.13%  '
    $ws.Range("E47").Value = '  -0.57%  '
    $ws.Range("D48").Value = '9.30'
    $ws.Range("E48").Value = '  +1.23%  '
    $ws.Range("D49").Value = '1.80'
    $ws.Range("E49").Value = '  +18.26%  '
    $ws.Range("D50").Value = '76.28'
    $ws.Range("E50").Value = '  +6.57%  '
    $ws.Range("D51").Value = '53.87'
    $ws.Range("E51").Value = '  +4.39%  '
